# G4 entfernt, I2C EEPROM hinzugefügt
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 35 was empty (the sheet data jumps from row 34 straight to row 48);
# just populate it with the new I2C EEPROM BOM line — no shifting of the
# rows below it.
$ws.Cells.Item(35, 1).Value = 1
$ws.Cells.Item(35, 2).Value = "EEPROM 8k SOT-23"
$ws.Cells.Item(35, 4).Value = "556-AT24CS08-STUM-T "

# Match the author's final selection.
$ws.Range("D35").Select()
